# Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) within specific bullet /
# impact paragraphs, per the commit's "Implement quantitative metrics
# highlighting across all resume formats" change.

$d = $word.ActiveDocument

# Word's Font.Color is a COLORREF-style long: 0x00BBGGRR (B<<16 | G<<8 | R).
# Target hex color is 2C3E50 (R=0x2C, G=0x3E, B=0x50).
$HighlightColor = (0x50 * 65536) + (0x3E * 256) + 0x2C

function Highlight-InParagraph($ParaText, $Terms) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        $rng = $p.Range
        $full = $rng.Text
        $trimmed = $full.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $ParaText) {
            $baseStart = $rng.Start
            $searchFrom = 0
            foreach ($term in $Terms) {
                $idx = $trimmed.IndexOf($term, $searchFrom)
                if ($idx -ge 0) {
                    $segStart = $baseStart + $idx
                    $segEnd = $segStart + $term.Length
                    $sub = $d.Range($segStart, $segEnd)
                    $sub.Font.Bold = 1
                    $sub.Font.Color = $HighlightColor
                    $searchFrom = $idx + $term.Length
                }
            }
            return $true
        }
    }
    return $false
}

$jobs = @(
    @('• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%', @('23%', '64%')),
    @('• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%', @('87%', '71%', '±4.2%', '±2.1%')),
    @('• Wrote RFP and analyzed bids from 1,200 vendors for research platform development', @('1,200')),
    @('• Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+', @('$400M', '$1B')),
    @('• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M', @('73.5%', '$4.7M')),
    @('• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%', @('87%', '71%'))
)

$failures = 0
foreach ($job in $jobs) {
    $ok = Highlight-InParagraph $job[0] $job[1]
    if (-not $ok) {
        $failures = $failures + 1
        Write-Host "NOT FOUND:" $job[0]
    }
}
Write-Host "Highlighted" ($jobs.Count - $failures) "of" $jobs.Count "paragraphs"
